$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Headers: BTec_Logo-Orange picture renamed image1.jpg -> image2.jpg
    for ($i = 1; $i -le 3; $i++) {
        $hf = $sec.Headers.Item($i)
        if ($hf.Exists) {
            foreach ($ishp in $hf.Range.InlineShapes) {
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $ishp.Name = "image2.jpg"
                }
            }
        }
    }

    # Footers: Pearson logo picture renamed image2.png -> image1.png
    for ($i = 1; $i -le 3; $i++) {
        $ft = $sec.Footers.Item($i)
        if ($ft.Exists) {
            foreach ($ishp in $ft.Range.InlineShapes) {
                if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $ishp.Name = "image1.png"
                }
            }
        }
    }
}
